{"js": "// Add a yellow highlight to the text \"(Partie forum UNDER PROGRESS)\"\n// (the parenthetical note marking that section of the design doc as\n// still under progress). The highlight is applied to the existing runs\n// that make up that phrase; nothing else in the paragraph (e.g. the\n// trailing page break) is touched.\nconst searchResults = context.document.body.search(\"(Partie forum UNDER PROGRESS)\", {\n  matchCase: true\n});\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find the text \"(Partie forum UNDER PROGRESS)\" in the document body.');\n}\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].font.highlightColor = \"Yellow\";\n}\n\nawait context.sync();\n", "ps1": "# Highlight the \"(Partie forum UNDER PROGRESS)\" note in yellow to flag\n# that section of the design doc as still a work in progress.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"(Partie forum UNDER PROGRESS)\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$found = $find.Execute()\n\nif ($found) {\n    # NOTE: setting HighlightColorIndex directly on the Range paints the\n    # whole enclosing paragraph (including the trailing page-break run),\n    # which is more than we want. Going through .Font restricts the\n    # highlight to exactly the matched run span.\n    $rng.Font.HighlightColorIndex = 7  # wdYellow\n}\n"}
